$d = $word.ActiveDocument

# --- Change 1 & 2: address paragraph - apply BodyText style and merge runs ---
$found = $d.Content.Find.Execute("<<Address_GLBL_Zip_Postal_Code_GLBL>> <<Address_GLBL_City_GLBL>>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found: $found"
